# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the bold/centered/bordered header style used by the rest of row 1 (e.g. column A1),
# then fill in the header text for the new columns AD, AE, AF.
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 30).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(1, 31).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(1, 32).PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item(1, 30).Value2 = "Wins"
$ws.Cells.Item(1, 31).Value2 = "Losses"
$ws.Cells.Item(1, 32).Value2 = "Ties"

# Data rows 2-41: every team has the same season record (80 wins, 82 losses, 0 ties).
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 80
    $ws.Cells.Item($r, 31).Value2 = 82
    $ws.Cells.Item($r, 32).Value2 = 0
}
